# Update thickness_nm (D), pore_size_nm (E) and recomputed flux_lmh (F)
# for the GO membrane simulation rows (rows 2-11) in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newThickness = 60
$newPoreSize  = 0.3

$fluxValues = @{
    2  = 0.3318117977528089
    3  = 0.6636235955056179
    4  = 0.9954353932584267
    5  = 1.327247191011236
    6  = 1.659058988764045
    7  = 1.990870786516853
    8  = 2.322682584269662
    9  = 2.654494382022472
    10 = 2.98630617977528
    11 = 3.318117977528089
}

foreach ($row in 2..11) {
    $ws.Cells.Item($row, 4).Value = $newThickness   # D: thickness_nm
    $ws.Cells.Item($row, 5).Value = $newPoreSize     # E: pore_size_nm
    $ws.Cells.Item($row, 6).Value = $fluxValues[$row] # F: flux_lmh
}

$wb.Save()
